$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9 (the "total" row), shifting the
# old row 9 (total) down to row 10 and the old row 10 (footer) down to row 11.
$ws.Rows.Item(9).Insert()

# Seed the new row 9 with the same formatting as row 8 (the previous data
# row) so it matches the existing data-row pattern, then fix up the row
# height which PasteSpecial(Formats) does not carry over.
$ws.Range("A8:Q8").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)
$ws.Rows.Item(9).RowHeight = 25.5

# Recreate the merges used by the data rows (A:B, C:G, H:K, L:M, N:O).
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# Populate the new product row.
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "QUICK NAIL  LOTION"
$ws.Range("H9").Value = "-23:0"

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "0"
$ws.Range("L9").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N9").Value = "85.00"

$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "1955.0000"
$ws.Range("P9").NumberFormat = "0.00"

$ws.Range("Q9").NumberFormat = "@"
$ws.Range("Q9").Value = "23:0"
$ws.Range("Q9").NumberFormat = "@"

# Update the total row (now shifted to row 10) with the new total.
$ws.Range("P10").Value = 2235
